# Applies the "Predicting Fingers" commit to the sock/finger-counting
# document:
#  A) Five bulleted paragraphs get their "<bullet>" run and "<tab>" run
#     swapped (tab now leads, bullet+space follow in the same run) so the
#     bullet glyph renders right before the list text instead of at the
#     left tab stop.
#  B) The split "signifi" / bookmark / "cantly" run sequence is merged
#     back into a single run reading "...significantly..." and the
#     _GoBack bookmark is dropped from that spot.
#  C) A new "Predicting Fingers" narrative (intro paragraph, a relocated
#     _GoBack bookmark paragraph, and four bulleted explanation
#     paragraphs) is appended at the end of the body, after the existing
#     trailing empty paragraph and before the sectPr.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$bullet = [char]0x2022

function Get-ParaIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13)
        if ($t -eq $needle) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# Part A: reorder "<bullet>" / "<tab>" runs in the five bulleted items.
# ---------------------------------------------------------------------

# Every one of the five bulleted paragraphs shares this pPr-holder
# rsid triplet in the source document; keep it so the edit stays minimal.
$bulletPAttrs = "w:rsidR='00101534' w:rsidRDefault='00101534' w:rsidP='00101534'"

# A1: the one paragraph that also carries <w:lastRenderedPageBreak/>.
$needle = "$bullet`tRandomly select one pair of same color socks."
$ix = Get-ParaIndex $d $needle
if ($ix -lt 0) { throw "Could not locate the 'Randomly select' paragraph" }
$xml = "<w:p $wNs $bulletPAttrs><w:r><w:lastRenderedPageBreak/><w:tab/><w:t xml:space='preserve'>$bullet </w:t></w:r><w:r><w:t>Randomly select one pair of same color socks.</w:t></w:r></w:p>"
$d.Paragraphs.Item($ix).Range.InsertXML($xml)

# A2-A5: the remaining plain "<bullet><tab>text" paragraphs.
$bulletTexts = @(
    "Selecting three pairs of same color socks.",
    "5 pairs of black socks (10 individual black socks-50% of the inventory)",
    "3 pairs of brown socks (6 individual brown socks-30% of the inventory)",
    "2 pairs of white socks (4 individual white socks-20% of the inventory)"
)
foreach ($txt in $bulletTexts) {
    $needle = "$bullet`t$txt"
    $ix = Get-ParaIndex $d $needle
    if ($ix -lt 0) { throw "Could not locate bulleted paragraph: $txt" }
    $xml = "<w:p $wNs $bulletPAttrs><w:r><w:tab/><w:t xml:space='preserve'>$bullet </w:t></w:r><w:r><w:t>$txt</w:t></w:r></w:p>"
    $d.Paragraphs.Item($ix).Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# Part B: merge the "signifi" / bookmark / "cantly..." runs, drop the
# bookmark that used to live inside that sentence.
# ---------------------------------------------------------------------

$tailNeedle = "Now for the other problem of selecting one matching pair of each color (total of three colors) lowers our rate significantly. I would recommend 3 pairs of socks (one pair of each color) leaving us with a random selection of 6 individual socks. I believe that there would be a 16% chance of success of randomly selecting 3 like pairs of matching socks."
$combinedNeedle = "`t$tailNeedle"
$ix = Get-ParaIndex $d $combinedNeedle
if ($ix -lt 0) {
    # Fall back to matching on the still-split wording, in case this runs
    # before any prior normalization.
    $splitNeedle = "`tNow for the other problem of selecting one matching pair of each color (total of three colors) lowers our rate signifi"
    $ix = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t.StartsWith($splitNeedle)) { $ix = $i; break }
    }
}
if ($ix -lt 0) { throw "Could not locate the 'Now for the other problem' paragraph" }
$xml = "<w:p $wNs w:rsidR='00F7796D' w:rsidRDefault='00F7796D' w:rsidP='00101534'><w:r><w:tab/><w:t>$tailNeedle</w:t></w:r></w:p>"
$d.Paragraphs.Item($ix).Range.InsertXML($xml)

# ---------------------------------------------------------------------
# Part C: append the new "Predicting Fingers" body paragraphs at the
# very end of the document, right before the sectPr.
# ---------------------------------------------------------------------

$lastIx = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIx)
if ($lastPara.Range.Text.TrimEnd([char]13) -ne "") {
    throw "Expected trailing empty paragraph after 'Predicting Fingers' heading"
}

$newParaXml = @(
    "<w:p $wNs><w:r><w:tab/></w:r><w:r><w:t>The problem here is developing a system of predicting an end count based on sequential system of the following:</w:t></w:r></w:p>",
    "<w:p $wNs><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>",
    "<w:p $wNs><w:r><w:tab/><w:t xml:space='preserve'>$bullet </w:t></w:r><w:r><w:t xml:space='preserve'>For the first five numbers, she always starts the numbering system with the thumb </w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>(number 1).</w:t></w:r></w:p>",
    "<w:p $wNs><w:r><w:tab/><w:t xml:space='preserve'>$bullet </w:t></w:r><w:r><w:t xml:space='preserve'>Followed by the first finger (number2), middle finger (number 3), ring finger (number </w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>4) ending with the little finger (number 5).</w:t></w:r></w:p>",
    "<w:p $wNs><w:r><w:tab/></w:r><w:r><w:t xml:space='preserve'>Then the next round of 5 numbers goes like this: </w:t></w:r></w:p>",
    "<w:p $wNs><w:r><w:tab/><w:t xml:space='preserve'>$bullet </w:t></w:r><w:r><w:t xml:space='preserve'>Ringer finger (number 6), middle finger (number 7), first finger (number 8), the thumb </w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>(number 9) and the first finger (number10).</w:t></w:r></w:p>"
)

foreach ($paraXml in $newParaXml) {
    $lastIx = $d.Paragraphs.Count
    $lastPara = $d.Paragraphs.Item($lastIx)
    $lastPara.Range.InsertParagraphBefore()
    $newIx = $lastIx
    $newPara = $d.Paragraphs.Item($newIx)
    $newPara.Range.InsertXML($paraXml)
}

Write-Host "Done."
